$d = $word.ActiveDocument

# The new UDP-vs-TCP sentence is inserted into the (empty) paragraph that
# holds the hidden "_GoBack" bookmark, immediately before the bookmark, as
# a sequence of separate <w:r> runs (matching how the source document
# stores them).
$bmRange = $d.Bookmarks.Item("_GoBack").Range
$bmPos = $bmRange.Start

# The texts below are the individual runs exactly as they appear in the
# target document (split the way the original author's edits happened to
# split them), in left-to-right reading order.
$runs = @(
    "In terms of UDP vs TCP, m",
    "ost games are built around the UDP, best effort communication, protocol. Compared to TC",
    "P, ",
    "UDP suffers",
    " ",
    "less from latency and is appropriate in situations where lost packets are not critical. UDP is a connection",
    " ",
    "less",
    " ",
    "protocol, so unlike TCP there is no guarantee of data-delivery. The most common use of UDP in games is to",
    " ",
    "inform each client on positions of each player where due to the high number of updates missing an update is",
    " ",
    "not critical",
    "."
)

# Inserting several runs back-to-back with InsertBefore at the same
# collapsed position normally lets Word coalesce them into a single run.
# Placing one throw-away character right after the bookmark first changes
# that: subsequent InsertBefore calls at the bookmark position then stay
# as discrete runs. Remove the throw-away marker again at the end.
$marker = "@@MARKER@@"
$anchor = $d.Range($bmPos, $bmPos)
$anchor.InsertAfter($marker)

for ($i = $runs.Length - 1; $i -ge 0; $i--) {
    $r = $d.Range($bmPos, $bmPos)
    $r.InsertBefore($runs[$i])
}

$found = $d.Content
[void]$found.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found.Find.Found) {
    $found.Delete()
}
